# Generate Report for Handoff
#
# The "b.md" file has been handed off again (new xlf packages generated),
# so its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff" across the Overview sheet and both per-locale
# sheets (zh-cn / de-de); the per-locale sheets also get the new handoff
# file name / datetime and an error detail message, and the "Error
# Detail" column is widened to fit it.

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"
$newDateTimeOverview = "2016-08-28 04:37:27"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c05a52e82a77d52ae13ad2a0776ca7523ff89079/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4bbc1e2c90c4c7a5998e981cd20764f5dd7daa1f/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $newDateTimeOverview

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 04:37:22"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1:P3").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" file.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $newDateTimeOverview
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1:P3").ColumnWidth = 39.166666666666664
